# Add two new columns (I: "I0", J: "IF") to the existing table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Copy the formatting of the last existing header cell (H1) onto the
# new header cells so they keep the same bold/bordered/centered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (rows 2-7) ---
$data = @{
    2 = @(1, 2)
    3 = @(1, 5)
    4 = @(1, 5)
    5 = @(7, 8)
    6 = @(4, 5)
    7 = @(6, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]   # column I
    $ws.Cells.Item($row, 10).Value = $vals[1]  # column J
}

$excel.CutCopyMode = 0
